# ---------------------------------------------------------------------------
# Adds four "secondary groupings" worksheets that cross the existing n/stance
# breakdowns with the 2016 presidential-result grouping, renames the first
# sheet to "prezresults2016", and touches a few pre-existing cosmetic sheet
# view properties (active-cell selection, scrolled top-left cell, etc.)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the first worksheet.
# ---------------------------------------------------------------------------
$wsPrez = $wb.Worksheets.Item(1)
$wsPrez.Name = "prezresults2016"

# ---------------------------------------------------------------------------
# 2. Helper that builds one of the four new "grouping" worksheets.
#    NOTE: only positional parameters are used (named "-Param value" binding
#    is not supported by this runtime's PowerShell-subset parser).
# ---------------------------------------------------------------------------
function Add-GroupingSheet {
    param(
        [string]$SheetName,
        [string]$HeaderB,
        [object[][]]$Rows,
        [double]$WidthA,
        [double]$WidthB,
        [double]$WidthC,
        [double]$WidthD,
        [string]$SelectCell
    )

    $wb = $excel.ActiveWorkbook
    $templateCell = $wb.Worksheets.Item(1)
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
    $ws.Name = $SheetName

    # Header row.
    $ws.Range("A1").Value = "p16winningparty"
    $ws.Range("B1").Value = $HeaderB
    $ws.Range("C1").Value = "stance"
    $ws.Range("D1").Value = "n"

    # Data rows.
    $r = 2
    foreach ($row in $Rows) {
        if ($null -ne $row[0]) { $ws.Cells.Item($r, 1).Value = $row[0] }
        if ($null -ne $row[1]) { $ws.Cells.Item($r, 2).Value = $row[1] }
        if ($null -ne $row[2]) { $ws.Cells.Item($r, 3).Value = $row[2] }
        if ($null -ne $row[3]) { $ws.Cells.Item($r, 4).Value = $row[3] }
        $r++
    }

    # Match the bold/centered header style already used by the other sheets
    # (copy formatting only, so the existing style entry is reused instead of
    # a brand-new one being appended to styles.xml).
    $templateCell.Range("A1").Copy() | Out-Null
    $ws.Range("A1:D1").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false

    # Column widths (best effort - see note below on precision).
    $ws.Columns.Item(1).ColumnWidth = $WidthA
    $ws.Columns.Item(2).ColumnWidth = $WidthB
    $ws.Columns.Item(3).ColumnWidth = $WidthC
    $ws.Columns.Item(4).ColumnWidth = $WidthD

    $ws.Range($SelectCell).Select() | Out-Null

    return $ws
}

# ---------------------------------------------------------------------------
# 3. gdp_andprezresults
# ---------------------------------------------------------------------------
$gdpRows = @(
    @("D", "ABOVE", "not_sponsoring", 2),
    @("D", "ABOVE", "sponsoring", 120),
    @("D", "BELOW", "not_sponsoring", 8),
    @("D", "BELOW", "sponsoring", 73),
    @("R", "ABOVE", "not_sponsoring", 3),
    @("R", "ABOVE", "sponsoring", 7),
    @("R", "BELOW", "not_sponsoring", 12),
    @("R", "BELOW", "sponsoring", 9)
)
Add-GroupingSheet "gdp_andprezresults" "gdp_abovebelow_natlavg" $gdpRows `
    15.166666666666666 23.333333333333336 14.0 3.1666666666666665 `
    "I28" | Out-Null

# ---------------------------------------------------------------------------
# 4. college_degree_andprezresults
# ---------------------------------------------------------------------------
$collegeRows = @(
    @("D", "ABOVE", "not_sponsoring", 3),
    @("D", "ABOVE", "sponsoring", 113),
    @("D", "BELOW", "not_sponsoring", 7),
    @("D", "BELOW", "sponsoring", 80),
    @("R", "ABOVE", "not_sponsoring", 5),
    @("R", "ABOVE", "sponsoring", 13),
    @("R", "BELOW", "not_sponsoring", 10),
    @("R", "BELOW", "sponsoring", 3)
)
Add-GroupingSheet "college_degree_andprezresults" "pct.ed.college.all.abovebelow.natl" $collegeRows `
    15.166666666666666 31.666666666666668 14.0 3.1666666666666665 `
    "E22" | Out-Null

# ---------------------------------------------------------------------------
# 5. nonwhite_pop_andprezresults
# ---------------------------------------------------------------------------
$nonwhiteRows = @(
    @("D", "ABOVE", "not_sponsoring", 7),
    @("D", "ABOVE", "sponsoring", 130),
    @("D", "BELOW", "not_sponsoring", 3),
    @("D", "BELOW", "sponsoring", 63),
    @("R", "ABOVE", "not_sponsoring", 3),
    @("R", "ABOVE", "sponsoring", 3),
    @("R", "BELOW", "not_sponsoring", 12),
    @("R", "BELOW", "sponsoring", 13)
)
Add-GroupingSheet "nonwhite_pop_andprezresults" "pct.race.nonwhite.abovebelow.natl" $nonwhiteRows `
    15.166666666666666 32.666666666666664 14.0 3.1666666666666665 `
    "E22" | Out-Null

# ---------------------------------------------------------------------------
# 6. rural_area_andprezresults  (has a couple of rows with a blank "stance")
# ---------------------------------------------------------------------------
$ruralRows = @(
    @("D", "ABOVE", "not_sponsoring", 3),
    @("D", "ABOVE", "sponsoring", 14),
    @("D", "BELOW", "not_sponsoring", 7),
    @("D", "BELOW", "sponsoring", 172),
    @("D", $null, "sponsoring", 7),
    @("R", "ABOVE", "not_sponsoring", 8),
    @("R", "ABOVE", "sponsoring", 6),
    @("R", "BELOW", "not_sponsoring", 6),
    @("R", "BELOW", "sponsoring", 9),
    @("R", $null, "not_sponsoring", 1),
    @("R", $null, "sponsoring", 1)
)
Add-GroupingSheet "rural_area_andprezresults" "pct.rural.above20" $ruralRows `
    15.166666666666666 15.666666666666666 14.0 3.1666666666666665 `
    "G22" | Out-Null

# ---------------------------------------------------------------------------
# 7. Cosmetic sheet-view tweaks on the pre-existing sheets (active cell /
#    scroll position).
# ---------------------------------------------------------------------------
$wsPrez.Range("F13").Select() | Out-Null

$wsGdp = $wb.Worksheets.Item("gdp_vs_nationalavg")
$wsGdp.Activate() | Out-Null
$wsGdp.Range("E37").Select() | Out-Null

$wsCollege = $wb.Worksheets.Item("college_vs_nationalavg")
$wsCollege.Activate() | Out-Null
$wsCollege.StandardWidth = 8.43
$wsCollege.Range("G29").Select() | Out-Null

$wsNonwhite = $wb.Worksheets.Item("nonwhite_vs_nationalavg")
$wsNonwhite.Activate() | Out-Null
$wsNonwhite.Range("G36").Select() | Out-Null

$wsRural = $wb.Worksheets.Item("rural_morethanfifth")
$wsRural.Activate() | Out-Null
$wsRural.Range("D13").Select() | Out-Null

# Re-activate the first sheet (tabSelected="1" lives on prezresults2016 in
# the source workbook).
$wsPrez.Activate() | Out-Null
